$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '24.869.51'
$ws.Range("E2").Value = '  +0.43%  '

$ws.Range("D3").Value = '1.710.72'
$ws.Range("E3").Value = '  +0.62%  '

$ws.Range("D4").Value = "'0.9989"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.40%  '

$ws.Range("D5").Value = "'317.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.05%  '

$ws.Range("D6").Value = "'0.9998"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.27%  '

$ws.Range("D7").Value = "'0.3930"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.39%  '

$ws.Range("D8").Value = "'0.4061"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.45%  '

$ws.Range("D9").Value = "'1.496"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.89%  '

$ws.Range("D10").Value = "'53.87"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.09%  '

$ws.Range("D11").Value = "'0.9977"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.56%  '

$ws.Range("D12").Value = "'0.08827"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.86%  '

$ws.Range("D13").Value = "'26.51"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +11.23%  '

$ws.Range("D14").Value = "'7.519"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.19%  '

$ws.Range("D15").Value = "'8.115"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.86%  '

$ws.Range("D16").Value = "'0.00001361"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.51%  '

$ws.Range("D17").Value = '1.725.89'
$ws.Range("E17").Value = '  +1.51%  '

$ws.Range("D18").Value = "'96.98"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.72%  '

$ws.Range("D19").Value = "'0.07240"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.57%  '

$ws.Range("D20").Value = "'20.88"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +4.55%  '

$ws.Range("D21").Value = "'7.290"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.12%  '

$ws.Range("D22").Value = "'0.9998"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.33%  '

$ws.Range("D23").Value = "'14.44"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.28%  '

$ws.Range("D24").Value = '24.852.19'
$ws.Range("E24").Value = '  +0.39%  '

$ws.Range("D25").Value = "'3.005"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -4.46%  '

$ws.Range("D26").Value = "'2.336"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.67%  '

$ws.Range("D27").Value = "'23.34"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.56%  '

$ws.Range("D28").Value = "'167.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.33%  '

$ws.Range("D29").Value = "'5.972"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +15.52%  '

$ws.Range("B30").Value = 'Filecoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D30").Value = "'8.600"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -7.87%  '

$ws.Range("B31").Value = 'BitcoinCash'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D31").Value = "'145.51"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +5.13%  '

$ws.Range("B32").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C32").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D32").Value = '1.916.60'
$ws.Range("E32").Value = '  +1.56%  '

$ws.Range("B33").Value = 'WEMIXTOKEN'
$ws.Range("C33").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D33").Value = "'2.241"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +14.43%  '

$ws.Range("D34").Value = "'0.08819"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.49%  '

$ws.Range("D35").Value = "'0.03161"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +4.33%  '

$ws.Range("D36").Value = "'1.051"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.08%  '

$ws.Range("D37").Value = "'7.239"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -9.52%  '

$ws.Range("D38").Value = "'0.2849"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.02%  '

$ws.Range("B39").Value = 'FraxShare'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D39").Value = "'10.96"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.11%  '

$ws.Range("B40").Value = 'TheSandbox'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D40").Value = "'0.8418"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +7.92%  '

$ws.Range("D41").Value = "'0.09246"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.46%  '

$ws.Range("D42").Value = "'14.17"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.98%  '

$ws.Range("D43").Value = "'1.478"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.48%  '

$ws.Range("E44").Value = '  +8.99%  '

$ws.Range("D45").Value = "'2.701"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.64%  '

$ws.Range("D46").Value = "'0.7449"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.88%  '

$ws.Range("D47").Value = "'4.254"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.33%  '

$ws.Range("D48").Value = "'1.403"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.90%  '

$ws.Range("D49").Value = "'1.000"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.18%  '

$ws.Range("D50").Value = "'141.49"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.93%  '

$ws.Range("D51").Value = "'0.08306"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.71%  '
